$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("L3").Value = "2017-02-22 08:25:08"
$wsZhCn.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/d71de91c25f3493ccb450101d23015e0cd4b8044/e2e/1d3eb6dd-6952-4a28-b9c4-63947b8ebd7c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/90894c2042f210f3f44d3e93013d54039f9baecd/e2e/1d3eb6dd-6952-4a28-b9c4-63947b8ebd7c.md."

$wsDeDe.Range("L3").Value = "2017-02-22 08:25:30"
$wsDeDe.Range("R3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/d71de91c25f3493ccb450101d23015e0cd4b8044/e2e/1d3eb6dd-6952-4a28-b9c4-63947b8ebd7c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test4/blob/90894c2042f210f3f44d3e93013d54039f9baecd/e2e/1d3eb6dd-6952-4a28-b9c4-63947b8ebd7c.md."
